$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the existing "Search tags..." text in C74 (non-privileged user
#    version). This reuses the same cell, so the shared string at the same
#    position gets replaced with the new wording.
# ---------------------------------------------------------------------------
$textNonPriv = @"
Search tags (which are optional) are taken into consideration if entered.
o Search only projects that I created
o Search others' projects
o Search for Products
o Search for Classes near zipcode ____
Notes:
--others' projects must be public
--products and classes must be active
--classes must start within 3 months
--will NEVER retrieve core projects
--when searching for a specific person's projects, user can enter that person's userName (email) as a search tag
"@
$ws.Range("C74").Value = $textNonPriv

# ---------------------------------------------------------------------------
# 2. Add the new "Search tags..." text for the privileged user row (C75),
#    plus the rest of row 75/76-80 content, in the same order the original
#    author entered them so new shared strings land at matching indexes.
# ---------------------------------------------------------------------------
$textPriv = @"
Search tags (which are optional) are taken into consideration if entered.
o Retrieve Core projects
o Search only projects that I created
o Search others' projects
o Search for Products
o Search for Classes near zipcode ____
Notes:
--retrieving core projects will ignore tags
--others' projects may be public or not
--products and classes may be active or not
--classes may start anytime--even in the past
--zipcode is optional when searching for classes
--will NEVER retrieve core projects
--when searching for a specific person's projects, user can enter that person's userName (email) as a search tag
"@
$ws.Range("C75").Value = $textPriv

$ws.Range("C76").Value = "Retrieve core projects"
$ws.Range("C77").Value = "Search only projects that I created"
$ws.Range("C78").Value = "Search others' products"
$ws.Range("C79").Value = "Search for Products"
$ws.Range("C80").Value = "Search for Classes"
$ws.Range("E76").Value = "Can save which will retain project id"
$ws.Range("D76").Value = "Can work on comic if authorized, changing but not adding; can work on System Types if authorized."

# ---------------------------------------------------------------------------
# 3. Adjust row heights to match the new wrapped content.
# ---------------------------------------------------------------------------
$ws.Rows.Item(74).RowHeight = 217.5
$ws.Rows.Item(75).RowHeight = 304.5
$ws.Rows.Item(76).RowHeight = 43.5

# ---------------------------------------------------------------------------
# 4. Update the view: scroll down and select E77, matching the saved
#    sheet view state from the edit.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 75
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E77").Select()
